{"js": "// Update the division-practice table: each populated row/column cell gets\n// a new \"a\u00f7b=\" expression per the target revision. Cells are addressed by\n// (row, col) on the first table so the two duplicate \"68\u00f77=\" entries are\n// each updated independently and correctly.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst replacements = [\n  { row: 0, col: 0, text: \"72\u00f79=\" },\n  { row: 0, col: 1, text: \"10\u00f78=\" },\n  { row: 0, col: 2, text: \"26\u00f78=\" },\n  { row: 0, col: 3, text: \"32\u00f72=\" },\n  { row: 0, col: 4, text: \"58\u00f76=\" },\n\n  { row: 4, col: 0, text: \"41\u00f75=\" },\n  { row: 4, col: 1, text: \"61\u00f75=\" },\n  { row: 4, col: 2, text: \"99\u00f74=\" },\n  { row: 4, col: 3, text: \"84\u00f74=\" },\n  { row: 4, col: 4, text: \"86\u00f79=\" },\n\n  { row: 8, col: 0, text: \"35\u00f74=\" },\n  { row: 8, col: 1, text: \"76\u00f72=\" },\n  { row: 8, col: 2, text: \"34\u00f75=\" },\n  { row: 8, col: 3, text: \"58\u00f74=\" },\n  { row: 8, col: 4, text: \"73\u00f76=\" },\n\n  { row: 12, col: 0, text: \"85\u00f76=\" },\n  { row: 12, col: 1, text: \"98\u00f73=\" },\n  { row: 12, col: 2, text: \"95\u00f77=\" },\n  { row: 12, col: 3, text: \"23\u00f76=\" },\n  { row: 12, col: 4, text: \"30\u00f73=\" },\n\n  { row: 16, col: 0, text: \"16\u00f76=\" },\n  { row: 16, col: 1, text: \"62\u00f79=\" },\n  { row: 16, col: 2, text: \"83\u00f75=\" },\n  { row: 16, col: 3, text: \"88\u00f72=\" },\n  { row: 16, col: 4, text: \"85\u00f77=\" },\n];\n\nfor (const r of replacements) {\n  const cell = table.getCell(r.row, r.col);\n  const paragraph = cell.body.paragraphs.getFirst();\n  const range = paragraph.getRange();\n  range.insertText(r.text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the division-practice table: each populated row/column cell gets\n# a new \"a\u00f7b=\" expression per the target revision. Cells are addressed by\n# (row, column) on the first table so the two duplicate \"68\u00f77=\" entries are\n# each updated independently and correctly.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{Row=1;  Col=1; Text=\"72\u00f79=\"},\n    @{Row=1;  Col=2; Text=\"10\u00f78=\"},\n    @{Row=1;  Col=3; Text=\"26\u00f78=\"},\n    @{Row=1;  Col=4; Text=\"32\u00f72=\"},\n    @{Row=1;  Col=5; Text=\"58\u00f76=\"},\n\n    @{Row=5;  Col=1; Text=\"41\u00f75=\"},\n    @{Row=5;  Col=2; Text=\"61\u00f75=\"},\n    @{Row=5;  Col=3; Text=\"99\u00f74=\"},\n    @{Row=5;  Col=4; Text=\"84\u00f74=\"},\n    @{Row=5;  Col=5; Text=\"86\u00f79=\"},\n\n    @{Row=9;  Col=1; Text=\"35\u00f74=\"},\n    @{Row=9;  Col=2; Text=\"76\u00f72=\"},\n    @{Row=9;  Col=3; Text=\"34\u00f75=\"},\n    @{Row=9;  Col=4; Text=\"58\u00f74=\"},\n    @{Row=9;  Col=5; Text=\"73\u00f76=\"},\n\n    @{Row=13; Col=1; Text=\"85\u00f76=\"},\n    @{Row=13; Col=2; Text=\"98\u00f73=\"},\n    @{Row=13; Col=3; Text=\"95\u00f77=\"},\n    @{Row=13; Col=4; Text=\"23\u00f76=\"},\n    @{Row=13; Col=5; Text=\"30\u00f73=\"},\n\n    @{Row=17; Col=1; Text=\"16\u00f76=\"},\n    @{Row=17; Col=2; Text=\"62\u00f79=\"},\n    @{Row=17; Col=3; Text=\"83\u00f75=\"},\n    @{Row=17; Col=4; Text=\"88\u00f72=\"},\n    @{Row=17; Col=5; Text=\"85\u00f77=\"}\n)\n\nforeach ($r in $replacements) {\n    $cell = $t.Cell($r.Row, $r.Col)\n    $rng = $cell.Range\n    $rng.End = $rng.End - 1\n    $rng.Text = $r.Text\n}\n"}
